# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 11:25"

# Row -> updated B:H values (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
$updates = @{
    4   = @(1745911, 108, 490151, 1153646, 0, 7,  102114)
    22  = @(57849,   257, 15572,  32889,   0, 24, 9388)
    44  = @(16628,   37,  15286,  674,     0, 23, 668)
    54  = @($null,   $null, 4827, 4712,    $null, $null, $null)
    62  = @(7629,    10,  6169,   1345,    $null, $null, $null)
    102 = @($null,   $null, 745,  714,     $null, $null, $null)
    109 = @(1076,    26,  823,    220,     0, 0,  33)
    110 = @(1067,    0,   1035,   28,      0, 0,  4)
    111 = @(1061,    4,   741,    296,     0, 1,  24)
    112 = @(1057,    $null, 779,  271,     $null, $null, 7)
    113 = @(1051,    $null, 929,  74,      $null, $null, 48)
    163 = @($null,   $null, 138,  1,       $null, $null, $null)
    189 = @(31,      3,   $null,  25,      $null, $null, $null)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals[0] -ne $null) { $ws.Cells.Item($row, 2).Value = $vals[0] }  # B
    if ($vals[1] -ne $null) { $ws.Cells.Item($row, 3).Value = $vals[1] }  # C
    if ($vals[2] -ne $null) { $ws.Cells.Item($row, 4).Value = $vals[2] }  # D
    if ($vals[3] -ne $null) { $ws.Cells.Item($row, 5).Value = $vals[3] }  # E
    if ($vals[4] -ne $null) { $ws.Cells.Item($row, 6).Value = $vals[4] }  # F
    if ($vals[5] -ne $null) { $ws.Cells.Item($row, 7).Value = $vals[5] }  # G
    if ($vals[6] -ne $null) { $ws.Cells.Item($row, 8).Value = $vals[6] }  # H
}
